$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new value looks numeric,
# so Excel stores them as text (matching the source data) instead of
# auto-converting to a number.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Update Price (D) and Volume(1h) (E) columns
$ws.Range("D2").Value = "57.489.25"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "3.089.88"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "516.64"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "141.79"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -1.09%  "
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("D10").Value = "0.109"
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("E11").Value = "  -1.36%  "
$ws.Range("D12").Value = "3.614.05"
$ws.Range("E12").Value = "  +1.22%  "
$ws.Range("E13").Value = "  +2.58%  "
$ws.Range("D14").Value = "25.66"
$ws.Range("E14").Value = "  -5.17%  "
$ws.Range("E15").Value = "  -1.85%  "
$ws.Range("D16").Value = "57.570.08"
$ws.Range("E16").Value = "  +0.73%  "
$ws.Range("D17").Value = "3.083.75"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").Value = "6.14"
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("D19").Value = "13.21"
$ws.Range("E19").Value = "  -1.54%  "
$ws.Range("D20").Value = "8.16"
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").Value = "334.84"
$ws.Range("E21").Value = "  +1.09%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").Value = "0.502"
$ws.Range("E23").Value = "  -1.22%  "
$ws.Range("D24").Value = "65.91"
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("E25").Value = "  +3.39%  "
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("D27").Value = "0.0₃0912"
$ws.Range("E27").Value = "  +2.20%  "
$ws.Range("D28").Value = "6.40"
$ws.Range("E28").Value = "  -5.26%  "
$ws.Range("D29").Value = "7.15"
$ws.Range("E29").Value = "  -0.93%  "
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("D31").Value = "20.86"
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("E32").Value = "  -3.66%  "
$ws.Range("D33").Value = "154.62"
$ws.Range("E33").Value = "  +2.65%  "
$ws.Range("D34").Value = "27.97"
$ws.Range("E34").Value = "  +10.65%  "
$ws.Range("E35").Value = "  -3.51%  "
$ws.Range("D36").Value = "5.89"
$ws.Range("E36").Value = "  -1.16%  "
$ws.Range("E37").Value = "  -0.89%  "
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("D39").Value = "3.128.46"
$ws.Range("E39").Value = "  +1.71%  "
$ws.Range("D40").Value = "36.77"
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("D41").Value = "3.88"
$ws.Range("E41").Value = "  -0.88%  "
$ws.Range("D42").Value = "0.672"
$ws.Range("E42").Value = "  +1.18%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").Value = "2.291.39"
$ws.Range("E44").Value = "  +3.94%  "
$ws.Range("D45").Value = "0.0256"
$ws.Range("E45").Value = "  +5.67%  "
$ws.Range("E46").Value = "  -1.44%  "
$ws.Range("E47").Value = "  -1.06%  "
$ws.Range("D48").Value = "20.05"
$ws.Range("E48").Value = "  -0.65%  "
$ws.Range("D49").Value = "5.90"
$ws.Range("E49").Value = "  -3.56%  "
$ws.Range("D50").Value = "252.84"
$ws.Range("E50").Value = "  +6.46%  "
$ws.Range("E51").Value = "  +0.99%  "
